# Auto-generated edit script: updates the cryptos price/volume table
# to match the refreshed data pulled on Fri Sep  6 18:33:09 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "53.746.02"
$ws.Range("E2").Value = "  -4.88%  "
# Row 3
$ws.Range("D3").Value = "2.234.71"
$ws.Range("E3").Value = "  -6.11%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.51"
$ws.Range("E5").Value = "  -3.15%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.36"
$ws.Range("E6").Value = "  -3.31%  "
# Row 7
$ws.Range("E7").Value = "  +0.24%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -5.02%  "
# Row 9
$ws.Range("D9").Value = "2.243.57"
$ws.Range("E9").Value = "  -6.01%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0919"
$ws.Range("E10").Value = "  -6.96%  "
# Row 11
$ws.Range("E11").Value = "  -1.23%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.74"
$ws.Range("E12").Value = "  +0.23%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.314"
$ws.Range("E13").Value = "  -3.57%  "
# Row 14
$ws.Range("D14").Value = "2.633.33"
$ws.Range("E14").Value = "  -6.08%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.12"
$ws.Range("E15").Value = "  -2.62%  "
# Row 16
$ws.Range("D16").Value = "53.653.40"
$ws.Range("E16").Value = "  -4.95%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000127"
$ws.Range("E17").Value = "  -3.87%  "
# Row 18
$ws.Range("D18").Value = "2.238.48"
$ws.Range("E18").Value = "  -4.92%  "
# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.58"
$ws.Range("E19").Value = "  -4.88%  "
# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.96"
$ws.Range("E20").Value = "  -1.82%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "298.90"
$ws.Range("E21").Value = "  -3.03%  "
# Row 22
$ws.Range("E22").Value = "  -2.19%  "
# Row 23
$ws.Range("E23").Value = "  +0.08%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.57"
$ws.Range("E24").Value = "  -3.01%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.62%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.363"
$ws.Range("E26").Value = "  -1.29%  "
# Row 27
$ws.Range("E27").Value = "  -3.88%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.01"
$ws.Range("E28").Value = "  -4.33%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.38"
$ws.Range("E29").Value = "  -1.69%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0685"
$ws.Range("E30").Value = "  -4.69%  "
# Row 31
$ws.Range("E31").Value = "  -3.12%  "
# Row 32
$ws.Range("E32").Value = "  -0.08%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.74"
$ws.Range("E33").Value = "  -0.66%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.13%  "
# Row 35
$ws.Range("E35").Value = "  -3.53%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.47"
$ws.Range("E36").Value = "  -0.84%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("E37").Value = "  -2.28%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  +4.74%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  -5.94%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.62"
$ws.Range("E40").Value = "  -1.40%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("E41").Value = "  -0.79%  "
# Row 42
$ws.Range("E42").Value = "  -2.85%  "
# Row 43
$ws.Range("E43").Value = "  -2.28%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.92"
$ws.Range("E44").Value = "  -6.25%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.64"
$ws.Range("E45").Value = "  -2.74%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0877"
$ws.Range("E46").Value = "  -3.43%  "
# Row 47
$ws.Range("E47").Value = "  -5.29%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "231.55"
$ws.Range("E48").Value = "  -4.30%  "
# Row 49
$ws.Range("E49").Value = "  -2.94%  "
# Row 50
$ws.Range("E50").Value = "  -3.60%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.06"
$ws.Range("E51").Value = "  -5.12%  "
